$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# randread_128k - IOPS (row 3)
$ws.Range("B3").Value = 2652
$ws.Range("C3").Value = 1850
$ws.Range("D3").Value = 3615
$ws.Range("E3").Value = 6330
$ws.Range("F3").Value = 11200
$ws.Range("G3").Value = 13100
# randread_128k - BW(mb/s) (row 4)
$ws.Range("B4").Value = 348.127232
$ws.Range("C4").Value = 242.221056
$ws.Range("D4").Value = 473.956352
$ws.Range("E4").Value = 829.423616
$ws.Range("F4").Value = 1473.24928
$ws.Range("G4").Value = 1718.616064
# randread_128k - lat_avg(us) (row 5)
$ws.Range("B5").Value = 375.67
$ws.Range("C5").Value = 1065.63
$ws.Range("D5").Value = 1087.76
$ws.Range("E5").Value = 1161.6
$ws.Range("F5").Value = 1357.14
$ws.Range("G5").Value = 2423.09
# randread_128k - p95_lat(us) (row 6)
$ws.Range("B6").Value = 668
$ws.Range("C6").Value = 1893
$ws.Range("D6").Value = 1958
$ws.Range("E6").Value = 2114
$ws.Range("G6").Value = 6128
# randread_128k - p99_lat(us) (row 7)
$ws.Range("B7").Value = 734
$ws.Range("C7").Value = 3130
$ws.Range("D7").Value = 2507
$ws.Range("E7").Value = 3392
$ws.Range("F7").Value = 3818
# randread_4k - IOPS (row 12)
$ws.Range("B12").Value = 3698
$ws.Range("C12").Value = 24600
$ws.Range("D12").Value = 51800
$ws.Range("E12").Value = 115000
$ws.Range("F12").Value = 58200
$ws.Range("G12").Value = 269000
# randread_4k - BW(mb/s) (row 13)
$ws.Range("B13").Value = 15.0994944
$ws.Range("C13").Value = 100.663296
$ws.Range("D13").Value = 211.812352
$ws.Range("E13").Value = 472.907776
$ws.Range("F13").Value = 238.026752
$ws.Range("G13").Value = 1103.101952
# randread_4k - lat_avg(us) (row 14)
$ws.Range("B14").Value = 269.13898
$ws.Range("C14").Value = 72.25129
$ws.Range("D14").Value = 65.063
$ws.Range("E14").Value = 53.11331
$ws.Range("F14").Value = 265.2566
$ws.Range("G14").Value = 111.89795
# randread_4k - p95_lat(us) (row 15)
$ws.Range("B15").Value = 1613.824
$ws.Range("C15").Value = 146.432
$ws.Range("D15").Value = 138.24
$ws.Range("E15").Value = 122.368
$ws.Range("F15").Value = 1564.672
$ws.Range("G15").Value = 325.632
# randread_4k - p99_lat(us) (row 16)
$ws.Range("B16").Value = 1777.664
$ws.Range("C16").Value = 197.632
$ws.Range("D16").Value = 175.104
$ws.Range("E16").Value = 166.912
$ws.Range("F16").Value = 1761.28
$ws.Range("G16").Value = 481.28
# randwrite_128k - IOPS (row 21)
$ws.Range("B21").Value = 2290
$ws.Range("C21").Value = 4511
$ws.Range("D21").Value = 9660
$ws.Range("E21").Value = 14100
$ws.Range("F21").Value = 16400
$ws.Range("G21").Value = 18500
# randwrite_128k - BW(mb/s) (row 22)
$ws.Range("B22").Value = 299.892736
$ws.Range("C22").Value = 591.3968640000001
$ws.Range("D22").Value = 1266.679808
$ws.Range("E22").Value = 1854.930944
$ws.Range("F22").Value = 2151.677952
$ws.Range("G22").Value = 2421.161984
# randwrite_128k - lat_avg(us) (row 23)
$ws.Range("B23").Value = 370.62
$ws.Range("C23").Value = 329.29
$ws.Range("D23").Value = 229.63
$ws.Range("E23").Value = 147.46
$ws.Range("F23").Value = 232.4
$ws.Range("G23").Value = 353.7
# randwrite_128k - p95_lat(us) (row 24)
$ws.Range("C24").Value = 457
$ws.Range("E24").Value = 375
$ws.Range("G24").Value = 486
# randwrite_128k - p99_lat(us) (row 25)
$ws.Range("B25").Value = 553
$ws.Range("C25").Value = 996
$ws.Range("D25").Value = 619
$ws.Range("E25").Value = 824
$ws.Range("F25").Value = 1037
$ws.Range("G25").Value = 963
# randwrite_4k - IOPS (row 30)
$ws.Range("B30").Value = 42300
$ws.Range("C30").Value = 131000
$ws.Range("D30").Value = 210000
$ws.Range("E30").Value = 303000
$ws.Range("F30").Value = 363000
$ws.Range("G30").Value = 421000
# randwrite_4k - BW(mb/s) (row 31)
$ws.Range("B31").Value = 173.01504
$ws.Range("C31").Value = 534.77376
$ws.Range("D31").Value = 861.929472
$ws.Range("E31").Value = 1241.513984
$ws.Range("F31").Value = 1487.929344
$ws.Range("G31").Value = 1723.858944
# randwrite_4k - lat_avg(us) (row 32)
$ws.Range("B32").Value = 21.2
$ws.Range("C32").Value = 11.4
$ws.Range("D32").Value = 12.12
$ws.Range("E32").Value = 12.8
$ws.Range("F32").Value = 19.77
$ws.Range("G32").Value = 28.64
# randwrite_4k - p95_lat(us) (row 33)
$ws.Range("B33").Value = 43
$ws.Range("C33").Value = 25
$ws.Range("D33").Value = 40
$ws.Range("E33").Value = 33
$ws.Range("F33").Value = 56
$ws.Range("G33").Value = 53
# randwrite_4k - p99_lat(us) (row 34)
$ws.Range("B34").Value = 55
$ws.Range("C34").Value = 42
$ws.Range("D34").Value = 51
$ws.Range("F34").Value = 92
$ws.Range("G34").Value = 122
# read_128k - IOPS (row 39)
$ws.Range("B39").Value = 10700
$ws.Range("C39").Value = 4394
$ws.Range("D39").Value = 7098
$ws.Range("E39").Value = 7160
$ws.Range("F39").Value = 7560
$ws.Range("G39").Value = 8114
# read_128k - BW(mb/s) (row 40)
$ws.Range("B40").Value = 1397.751808
$ws.Range("C40").Value = 575.668224
$ws.Range("D40").Value = 930.086912
$ws.Range("E40").Value = 938.47552
$ws.Range("F40").Value = 990.90432
$ws.Range("G40").Value = 1063.256064
# read_128k - lat_avg(us) (row 41)
$ws.Range("B41").Value = 92.47
$ws.Range("C41").Value = 412.15
$ws.Range("D41").Value = 548.05
$ws.Range("E41").Value = 1018.02
$ws.Range("F41").Value = 1984.3
$ws.Range("G41").Value = 3733.28
# read_128k - p95_lat(us) (row 42)
$ws.Range("B42").Value = 285
$ws.Range("C42").Value = 1827
$ws.Range("D42").Value = 1762
$ws.Range("E42").Value = 2999
$ws.Range("F42").Value = 5407
$ws.Range("G42").Value = 10159
# read_128k - p99_lat(us) (row 43)
$ws.Range("B43").Value = 379
$ws.Range("C43").Value = 2999
$ws.Range("D43").Value = 2024
$ws.Range("E43").Value = 4146
$ws.Range("F43").Value = 6849
$ws.Range("G43").Value = 13173
# read_4k - IOPS (row 48)
$ws.Range("B48").Value = 131000
$ws.Range("C48").Value = 158000
$ws.Range("D48").Value = 207000
$ws.Range("E48").Value = 246000
$ws.Range("F48").Value = 246000
$ws.Range("G48").Value = 246000
# read_4k - BW(mb/s) (row 49)
$ws.Range("B49").Value = 534.77376
$ws.Range("C49").Value = 644.87424
$ws.Range("D49").Value = 848.297984
$ws.Range("E49").Value = 1006.63296
$ws.Range("F49").Value = 1009.778688
$ws.Range("G49").Value = 1005.584384
# read_4k - lat_avg(us) (row 50)
$ws.Range("B50").Value = 7.14975
$ws.Range("C50").Value = 10.04875
$ws.Range("D50").Value = 17.99282
$ws.Range("E50").Value = 29.62835
$ws.Range("F50").Value = 60.65716
$ws.Range("G50").Value = 124.11287
# read_4k - p95_lat(us) (row 51)
$ws.Range("B51").Value = 2.16
$ws.Range("C51").Value = 2.224
$ws.Range("D51").Value = 2.768
$ws.Range("E51").Value = 3.056
$ws.Range("F51").Value = 2.992
$ws.Range("G51").Value = 2.992
# read_4k - p99_lat(us) (row 52)
$ws.Range("B52").Value = 144.384
$ws.Range("C52").Value = 162.816
$ws.Range("D52").Value = 514.048
$ws.Range("E52").Value = 1597.44
$ws.Range("F52").Value = 3227.648
$ws.Range("G52").Value = 7241.728
# write_128k - IOPS (row 57)
$ws.Range("B57").Value = 2503
$ws.Range("C57").Value = 4452
$ws.Range("D57").Value = 8827
$ws.Range("E57").Value = 11500
$ws.Range("F57").Value = 12300
$ws.Range("G57").Value = 12400
# write_128k - BW(mb/s) (row 58)
$ws.Range("B58").Value = 328.204288
$ws.Range("C58").Value = 584.056832
$ws.Range("D58").Value = 1156.579328
$ws.Range("E58").Value = 1509.94944
$ws.Range("F58").Value = 1609.56416
$ws.Range("G58").Value = 1625.2928
# write_128k - lat_avg(us) (row 59)
$ws.Range("B59").Value = 316.83
$ws.Range("C59").Value = 288.98
$ws.Range("D59").Value = 179.2
$ws.Range("E59").Value = 151.2
$ws.Range("F59").Value = 219.69
$ws.Range("G59").Value = 352.64
# write_128k - p95_lat(us) (row 60)
$ws.Range("B60").Value = 474
$ws.Range("C60").Value = 400
$ws.Range("D60").Value = 343
$ws.Range("E60").Value = 367
$ws.Range("F60").Value = 445
$ws.Range("G60").Value = 469
# write_128k - p99_lat(us) (row 61)
$ws.Range("B61").Value = 519
$ws.Range("C61").Value = 457
$ws.Range("D61").Value = 424
$ws.Range("E61").Value = 424
$ws.Range("F61").Value = 502
$ws.Range("G61").Value = 668
# write_4k - IOPS (row 66)
$ws.Range("B66").Value = 75200
$ws.Range("C66").Value = 151000
$ws.Range("D66").Value = 224000
$ws.Range("E66").Value = 284000
$ws.Range("F66").Value = 307000
$ws.Range("G66").Value = 348000
# write_4k - BW(mb/s) (row 67)
$ws.Range("B67").Value = 308.281344
$ws.Range("C67").Value = 617.611264
$ws.Range("D67").Value = 917.504
$ws.Range("E67").Value = 1162.870784
$ws.Range("F67").Value = 1256.194048
$ws.Range("G67").Value = 1427.111936
# write_4k - lat_avg(us) (row 68)
$ws.Range("B68").Value = 10.52
$ws.Range("E68").Value = 10.72
$ws.Range("F68").Value = 18.08
$ws.Range("G68").Value = 27.93
# write_4k - p95_lat(us) (row 69)
$ws.Range("B69").Value = 31
$ws.Range("C69").Value = 21
$ws.Range("D69").Value = 16
$ws.Range("F69").Value = 48
$ws.Range("G69").Value = 49
# write_4k - p99_lat(us) (row 70)
$ws.Range("B70").Value = 40
$ws.Range("D70").Value = 39
$ws.Range("E70").Value = 46
$ws.Range("F70").Value = 77
$ws.Range("G70").Value = 121
